$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..153 (A: angle in radians, B: =SIN(A)).
# Row 153 (A153 = 1.52) was copied/filled down through row 169, repeating the
# same constant in column A and extending the SIN() formula in column B.
for ($r = 154; $r -le 169; $r++) {
    $ws.Range("A$r").Value2 = 1.52
}

# Assign the formula across the whole new block in one shot so it is stored
# as a (shared) formula, like a fill-down/copy-down would produce.
$ws.Range("B154:B169").Formula = "=SIN(A154)"

# Restore the view/selection state: the user scrolled down and ended up with
# E160 selected.
$ws.Range("E160").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 88
    $win.ScrollColumn = 1
} catch {
    # View-scrolling is best-effort; ignore if unsupported.
}
